$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.25
$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.19
$ws.Range("P10").Value = 4.65
$ws.Range("S10").Value = 1.24
$ws.Range("T10").Value = 3.65
$ws.Range("G14").Value = 3.6
$ws.Range("I14").Value = 1.8
$ws.Range("J14").Value = 3.75
$ws.Range("L14").Value = 2.3
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 23
$ws.Range("Q14").Value = 1.4
$ws.Range("R14").Value = 2.88
$ws.Range("X14").Value = 23
$ws.Range("AA14").Value = 23
$ws.Range("AI14").Value = 13
$ws.Range("AL14").Value = 13
$ws.Range("AP14").Value = 19
$ws.Range("AS14").Value = 81
$ws.Range("AW14").Value = 4.5
$ws.Range("AZ14").Value = 26
$ws.Range("BA14").Value = 41
$ws.Range("BC14").Value = 201
$ws.Range("Q15").Value = 1.67
$ws.Range("Q16").Value = 1.83
$ws.Range("R16").Value = 2.03
$ws.Range("Q17").Value = 1.44
$ws.Range("Q19").Value = 1.7
$ws.Range("Q20").Value = 1.5
$ws.Range("R20").Value = 2.5
$ws.Range("U20").Value = 1.67
$ws.Range("U21").Value = 1.53
$ws.Range("V21").Value = 2.38
$ws.Range("U22").Value = 1.57
$ws.Range("V27").Value = 1.69
$ws.Range("U28").Value = 1.69
$ws.Range("G29").Value = 3.75
$ws.Range("H29").Value = 3.6
$ws.Range("I29").Value = 1.87
$ws.Range("J29").Value = 4.1
$ws.Range("K29").Value = 2.18
$ws.Range("L29").Value = 2.47
$ws.Range("T29").Value = 2.82
$ws.Range("W29").Value = 12
$ws.Range("X29").Value = 21
$ws.Range("Y29").Value = 12.5
$ws.Range("Z29").Value = 55
$ws.Range("AH29").Value = 7.8
$ws.Range("AI29").Value = 9.25
$ws.Range("AK29").Value = 16
$ws.Range("AN29").Value = 5.6
$ws.Range("AS29").Value = 300
$ws.Range("AT29").Value = 2.82
$ws.Range("AW29").Value = 3.8
$ws.Range("AX29").Value = 9.5
$ws.Range("AY29").Value = 18
$ws.Range("BB29").Value = 250
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("G32").Value = 2.38
$ws.Range("I32").Value = 2.7
$ws.Range("L32").Value = 3.2
$ws.Range("Q32").Value = 1.57
$ws.Range("U32").Value = 1.5
$ws.Range("W32").Value = 12
$ws.Range("X32").Value = 15
$ws.Range("AJ32").Value = 10
$ws.Range("AL32").Value = 19
$ws.Range("R33").Value = 1.53
$ws.Range("V33").Value = 1.75
$ws.Range("U34").Value = 1.91
$ws.Range("V34").Value = 1.91
$ws.Range("Q36").Value = 1.98
$ws.Range("R36").Value = 1.83
$ws.Range("U36").Value = 1.69
$ws.Range("Q39").Value = 1.93
$ws.Range("R39").Value = 1.93
